$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used (B,C,D,E,F,I,J,K,L,M,N) mapped to column indices
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(1.02, 1.03830827859863, 1.041893320539708, 1.046622765366318, 1.056740132583072, 1.040684779555896, 1.043406421329438, 1.044671547752643, 1.04938768418819, 1.059477015119916, 1.044888178128115),
    @(1.02, 1.039225788613012, 1.04259793952127, 1.047511494830233, 1.057878063177577, 1.040935700576638, 1.043968830998331, 1.045187165651309, 1.050087893307484, 1.060427819577611, 1.045451386483247),
    @(1.02, 1.039819702537381, 1.043053930912896, 1.048087199874871, 1.058615534116578, 1.041096773343455, 1.044332338092738, 1.045520186858627, 1.050540976741385, 1.061043620399547, 1.045815409799423),
    @(1.02, 1.040069436627061, 1.043245641828417, 1.048329377910945, 1.058925842360395, 1.041164178995864, 1.044485057633046, 1.045660040154958, 1.050731452525985, 1.061302637986782, 1.045968346218996),
    @(1.02, 1.040111371143819, 1.043277831610696, 1.04837004953118, 1.058977960611444, 1.041175478546308, 1.044510694096741, 1.045683513387067, 1.050763434196194, 1.061346136069233, 1.04599401908941),
    @(1.02, 1.03982303929059, 1.043056492517389, 1.048090435275095, 1.058619679389074, 1.041097675237234, 1.044334379128082, 1.04552205617105, 1.050543521891077, 1.061047080876021, 1.045817453733271),
    @(1.02, 1.038618308594031, 1.042131437536494, 1.046922983005476, 1.057124462210958, 1.040769846254991, 1.043596574638731, 1.044845930692502, 1.04962432238232, 1.059798226127876, 1.045078601476918),
    @(1.02, 1.036497169163338, 1.040501853621195, 1.044870711566019, 1.054498567896851, 1.040182313909408, 1.042293366854962, 1.043649815371922, 1.048004620912013, 1.05760195874566, 1.043773542988616),
    @(1.02, 1.035084301804099, 1.039415859531166, 1.043505903708899, 1.052753977631966, 1.03978403011017, 1.041422521997335, 1.042849297872576, 1.046924897197846, 1.056140765137475, 1.042901461431404),
    @(1.02, 1.034472815147954, 1.038945720665022, 1.042915739017801, 1.051999983709918, 1.03961001074408, 1.041044962144885, 1.042501938185608, 1.046457392703919, 1.055508768779807, 1.042523365400657),
    @(1.02, 1.034245726707986, 1.038771106903581, 1.042696647616227, 1.051720131520344, 1.039545138317939, 1.04090464845757, 1.042372804389153, 1.046283744888805, 1.055274124517541, 1.042382852451812),
    @(1.02, 1.03429443590435, 1.03880856137173, 1.042743637927966, 1.051780151038428, 1.03955906424989, 1.040934749427974, 1.042400508944189, 1.046320992734441, 1.055324451641373, 1.042412996169046),
    @(1.02, 1.034454043025887, 1.038931286688703, 1.04289762635919, 1.051976846656772, 1.039604653137595, 1.041033365218997, 1.042491266163587, 1.046443038832156, 1.055489370818797, 1.042511752005804),
    @(1.02, 1.034552388216658, 1.039006904060676, 1.042992519913547, 1.052098065790396, 1.039632710977783, 1.041094116256109, 1.042547170269243, 1.046518236034893, 1.055590997207727, 1.042572589316356),
    @(1.02, 1.035124890420586, 1.03944706338415, 1.043545088059252, 1.05280404780395, 1.039795546377724, 1.04144756938494, 1.042872335670497, 1.046955924476905, 1.056182723669977, 1.042926544389172),
    @(1.02, 1.035484085213649, 1.039723192293214, 1.043891915936294, 1.053247273656974, 1.039897271432983, 1.041669153645271, 1.043076108357918, 1.047230481534237, 1.05655408861393, 1.043148443324569),
    @(1.02, 1.03569362579986, 1.039884263545961, 1.044094292485195, 1.053505937265937, 1.039956455402248, 1.041798353800495, 1.043194895000587, 1.047390628191759, 1.056770767975768, 1.043277826958833),
    @(1.02, 1.035765078509765, 1.039939186329969, 1.044163310812159, 1.0535941581706, 1.039976610058944, 1.041842399857569, 1.043235386191124, 1.047445234409972, 1.056844661633535, 1.04332193556636),
    @(1.02, 1.035445544075053, 1.039693565238821, 1.043854696542122, 1.053199705494205, 1.039886372864372, 1.041645384524189, 1.043054252765414, 1.0472010239306, 1.056514237582859, 1.043124640448609),
    @(1.02, 1.034407041460554, 1.03889514667011, 1.042852277227814, 1.051918918767318, 1.039591234799079, 1.041004327284787, 1.042464543425614, 1.046407099177452, 1.055440803263723, 1.042482672834397),
    @(1.02, 1.033754353964669, 1.038393246689398, 1.042222722977529, 1.051114879953418, 1.039404316825999, 1.040600857870545, 1.042093139837004, 1.045907951653196, 1.054766513768606, 1.042078630447307),
    @(1.02, 1.034100331037531, 1.038659303740101, 1.042556394312067, 1.051540998149656, 1.039503533669289, 1.040814783388884, 1.042290087355002, 1.046172556640711, 1.055123908285162, 1.042292859764421),
    @(1.02, 1.035462959069047, 1.039706952398835, 1.043871514144185, 1.053221199074642, 1.039891297923552, 1.041656124909649, 1.043064128584845, 1.047214334547018, 1.056532244336859, 1.043135396086648),
    @(1.02, 1.037045322127676, 1.040923075710839, 1.045400682744807, 1.055176369071516, 1.04033536988456, 1.042630640175763, 1.043959591188064, 1.048423342369834, 1.058169223653046, 1.044111295276219)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowIndex, $cols[$j]).Value = $rowValues[$j]
    }
}
